# Generate Report for Handback
# Updates the localization-status workbook:
#  - Overview sheet: status text "In Translation" -> "Handed back: in sync with en-US"
#    for the two localized-language columns (E/F), and widen those columns.
#  - zh-cn / de-de sheets: fill in "Latest Target File" (I) and
#    "Latest Handback File" (J) for both rows, add hyperlinks on the new
#    target-file cells (mirroring the existing source-file hyperlinks in A),
#    widen the Status/Target/Handback columns, and stamp the handback
#    timestamp into "Latest Handback DateTime" (K).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Overview sheet: status text + column widths
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"

$wsOverview.Columns("E").ColumnWidth = 29.1
$wsOverview.Columns("F").ColumnWidth = 29.1

# ---------------------------------------------------------------------------
# Helper data shared by the two localized-language sheets
# ---------------------------------------------------------------------------
$sourceMd1 = "07271aef-b338-4d0e-bb0c-0439fea621d8.md"
$sourceMd2 = "241b475f-77c1-4ae8-972a-0a5fd4b7fdd4.md"
$sourceUrl1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3e6cbbdb89f9bc3dba1b43357a5a9e27f88950fc/e2e/07271aef-b338-4d0e-bb0c-0439fea621d8.md"
$sourceUrl2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3e6cbbdb89f9bc3dba1b43357a5a9e27f88950fc/e2e/241b475f-77c1-4ae8-972a-0a5fd4b7fdd4.md"

# ---------------------------------------------------------------------------
# 2. zh-cn sheet
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("C2").Value = "Handed back: in sync with en-US"
$wsZh.Range("C3").Value = "Handed back: in sync with en-US"

$wsZh.Range("I2").Value = $sourceMd1
$wsZh.Range("J2").Value = "07271aef-b338-4d0e-bb0c-0439fea621d8.8ae8057a783306840244472b8dc505b829da902b.zh-cn.xlf"
$wsZh.Range("I3").Value = $sourceMd2
$wsZh.Range("J3").Value = "241b475f-77c1-4ae8-972a-0a5fd4b7fdd4.da02482e99af397be193255abd4eb61f28d0acaf.zh-cn.xlf"

$wsZh.Range("K2").Value = "2016-09-02 20:29:28"
$wsZh.Range("K3").Value = "2016-09-02 20:29:28"

# Re-create the hyperlinks so the new Latest-Target-File cells (I2/I3) get
# their own link, in the same order Excel would assign relationship ids:
# A2, I2, A3, I3.
$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $sourceUrl1, "", "", $sourceMd1)
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $sourceUrl1, "", "", $sourceMd1)
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $sourceUrl2, "", "", $sourceMd2)
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), $sourceUrl2, "", "", $sourceMd2)

$wsZh.Columns("C").ColumnWidth = 29.1
$wsZh.Columns("I").ColumnWidth = 39.17
$wsZh.Columns("J").ColumnWidth = 39.17

# ---------------------------------------------------------------------------
# 3. de-de sheet
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C2").Value = "Handed back: in sync with en-US"
$wsDe.Range("C3").Value = "Handed back: in sync with en-US"

$wsDe.Range("I2").Value = $sourceMd1
$wsDe.Range("J2").Value = "07271aef-b338-4d0e-bb0c-0439fea621d8.8ae8057a783306840244472b8dc505b829da902b.de-de.xlf"
$wsDe.Range("I3").Value = $sourceMd2
$wsDe.Range("J3").Value = "241b475f-77c1-4ae8-972a-0a5fd4b7fdd4.da02482e99af397be193255abd4eb61f28d0acaf.de-de.xlf"

$wsDe.Range("K2").Value = "2016-09-02 20:29:36"
$wsDe.Range("K3").Value = "2016-09-02 20:29:36"

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $sourceUrl1, "", "", $sourceMd1)
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $sourceUrl1, "", "", $sourceMd1)
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $sourceUrl2, "", "", $sourceMd2)
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), $sourceUrl2, "", "", $sourceMd2)

$wsDe.Columns("C").ColumnWidth = 29.1
$wsDe.Columns("I").ColumnWidth = 39.17
$wsDe.Columns("J").ColumnWidth = 39.17
